# Apply the Alvearie FHIR IG regeneration update:
#  - Metadata sheet: bump Version, Date, Publisher, replace duplicate
#    "Contact" row with a single "Jurisdiction" row, shifting later rows up.
#  - Elements sheet: refresh root Extension's Short/Definition text to match
#    the new Title/Description.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact / No display for ContactDetail" row (row 11),
# which shifts every following row up by one.
$meta.Rows.Item(11).Delete()

# Row 10 (was the first "Contact" row) becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short (K) / Definition (L) now mirror the new Title/Description
$elements.Range("K2").Value = "Episode Allowed Amount Outpatient (USD)"
$elements.Range("L2").Value = "Allowed amount for an out-patient episode, in USD"
